$wb = $excel.ActiveWorkbook

$wsInsumos      = $wb.Worksheets.Item("Insumos")
$wsProcesados   = $wb.Worksheets.Item("Procesados")
$wsNoProcesados = $wb.Worksheets.Item("NoProcesados")

# --- Insumos: the duplicate "Iván García" record stays, and the previously
#     blank rows 4-5 are filled in with the validated "Allisson Flores
#     Espinoza" record (a repeated entry, caught by the validation pass) ---
$wsInsumos.Range("A4").Value = "Allisson"
$wsInsumos.Range("B4").Value = "Flores"
$wsInsumos.Range("C4").Value = "Espinoza"

$wsInsumos.Range("A5").Value = "Allisson"
$wsInsumos.Range("B5").Value = "Flores"
$wsInsumos.Range("C5").Value = "Espinoza"

# --- Procesados: row 3 actually corresponds to the same ivan.garcia email
#     (correcting the earlier mistaken ivan.hernandez address), and the new
#     validated records now produce their own generated emails ---
$wsProcesados.Range("A3").Value  = "ivan.garcia@beeckerco.com"
$wsProcesados.Range("A4").Value  = "allisson.flores@beeckerco.com"
$wsProcesados.Range("A5").Value  = "allisson.flores@beeckerco.com"
$wsProcesados.Range("A6").Value  = "eunice.@beeckerco.com"

# --- NoProcesados: everything got processed this run, so the previously
#     flagged rows are cleared out ---
$wsNoProcesados.Range("A2:D3").ClearContents()

# --- Restore each sheet's last-used selection, then land on "Procesados"
#     as the active sheet/tab, matching the saved workbook UI state ---
$wsNoProcesados.Activate()
$wsNoProcesados.Range("A2:E9").Select()

$wsProcesados.Activate()
$wsProcesados.Range("A2:C11").Select()
